$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A; this shifts the existing A:D data to B:E
$ws.Columns("A").Insert()

# Update header text to add spaces around "=" in the "n=" counts
$ws.Range("D1").Value = "Treatment at T2 (n = 5080)"
$ws.Range("E1").Value = "Control at T2 (n = 745)"

# Update category labels to add spaces around "=" in the "P=" values
$ws.Range("B3").Value = "Gender (P = 0.006)"
$ws.Range("B10").Value = "Interested in News (P = 0.000)"

# Give column A (rows 2-15) the same formatting (bold, border, center/top
# alignment) as the header row, matching the rest of the table's styling,
# without touching the (already empty) values in those cells.
$ws.Range("B1").Copy()
$ws.Range("A2:A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
